$d = $word.ActiveDocument

# 1. Table cell color rename: "Giallo" -> "Mandarino"
$d.Content.Find.Execute("Giallo", $true, $false, $false, $false, $false, $true, 1, $false, "Mandarino", 2)

# 2. Append an underlined space run right after the "COSE DA FARE" heading text
$pCose = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "COSE DA FARE") {
        $pCose = $cand
        break
    }
}
$insertPoint = $d.Range($pCose.Range.End - 1, $pCose.Range.End - 1)
$insertPoint.InsertAfter(" ")
$insertPoint.Font.Underline = 1

# 3. Find the "Mockup Mobile" and "Mockup Web" heading paragraphs (Titolo2 style) and
#    the bullet paragraphs that belong to each of them, then collapse the whole
#    "COSE DA FARE" sub-section down to a single "Mockup Web" heading followed by the
#    last bullet ("Finire di farli con i dovuti accorgimenti").
$pMobileIdx = -1
$pWebIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Mockup Mobile") { $pMobileIdx = $i }
    if ($t -eq "Mockup Web") { $pWebIdx = $i }
}

# Delete everything between the "Mockup Mobile" heading (exclusive) and the
# "Mockup Web" heading (exclusive) -- i.e. all the Mockup Mobile bullet points.
$startP = $d.Paragraphs.Item($pMobileIdx + 1)
$endP = $d.Paragraphs.Item($pWebIdx - 1)
$r = $d.Range($startP.Range.Start, $endP.Range.End)
$r.Delete()

# Rename the "Mockup Mobile" heading to "Mockup Web" (it keeps the Titolo2 style/props).
$pMobile = $d.Paragraphs.Item($pMobileIdx)
$pMobile.Range.Text = "Mockup Web"

# Recompute indices: the old "Mockup Web" heading is now immediately after the
# renamed heading. Delete that old heading plus its first two bullets, keeping
# the final bullet ("Finire di farli con i dovuti accorgimenti").
$oldWebHeadingIdx = $pMobileIdx + 1
$startP2 = $d.Paragraphs.Item($oldWebHeadingIdx)
$endP2 = $d.Paragraphs.Item($oldWebHeadingIdx + 2)
$r2 = $d.Range($startP2.Range.Start, $endP2.Range.End)
$r2.Delete()

Write-Output "Edit complete"
